# The document's heading hierarchy is being flattened: every sub-heading
# that currently uses Heading 3 / Heading 4 / Heading 5 is promoted to
# Heading 2 (this is how the commit "re-enables" the generated build,
# which only ever emits Heading 1 / Heading 2 paragraph styles).
#
# (The accompanying bookmark-id renumbering and hyperlink relationship-id
# shuffling visible in the source diff are artifacts of the document
# generator re-running - the bookmark targets/names and hyperlink URLs
# themselves are unchanged, so no user-visible/content edit is needed for
# those beyond what naturally happens when the package is rebuilt.)

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $styleName = $para.Style.NameLocal
    if ($styleName -eq "Heading 3" -or $styleName -eq "Heading 4" -or $styleName -eq "Heading 5") {
        $para.Style = "Heading 2"
    }
}
